$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with the latest scraped values.
# Numeric-looking Price values are prefixed with a leading apostrophe so Excel
# keeps them as literal text (matching the original inlineStr/text storage)
# instead of re-interpreting them as numbers.
$ws.Range("D2").Value = '29.238.07'
$ws.Range("D3").Value = '1.863.44'
$ws.Range("E3").Value = '  +0.80%  '
$ws.Range("D4").Value = '''0.9998'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '''0.7080'
$ws.Range("E5").Value = '  +0.70%  '
$ws.Range("D6").Value = '''237.92'
$ws.Range("E6").Value = '  -0.35%  '
$ws.Range("D7").Value = '''1.0000'
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '''0.08185'
$ws.Range("E8").Value = '  +9.20%  '
$ws.Range("D9").Value = '''0.3038'
$ws.Range("E9").Value = '  -0.59%  '
$ws.Range("D10").Value = '''23.33'
$ws.Range("E10").Value = '  -0.40%  '
$ws.Range("D11").Value = '''0.08172'
$ws.Range("E11").Value = '  +0.53%  '
$ws.Range("D12").Value = '1.831.93'
$ws.Range("E12").Value = '  -1.11%  '
$ws.Range("D13").Value = '''5.172'
$ws.Range("E13").Value = '  -1.09%  '
$ws.Range("D14").Value = '''0.7088'
$ws.Range("D15").Value = '''89.29'
$ws.Range("E15").Value = '  +0.65%  '
$ws.Range("D16").Value = '29.234.37'
$ws.Range("E16").Value = '  +0.09%  '
$ws.Range("D17").Value = '''0.000007905'
$ws.Range("E17").Value = '  +3.59%  '
$ws.Range("D18").Value = '''5.792'
$ws.Range("E18").Value = '  +0.40%  '
$ws.Range("D19").Value = '''13.39'
$ws.Range("E19").Value = '  +2.52%  '
$ws.Range("D20").Value = '''236.77'
$ws.Range("D21").Value = '''0.9993'
$ws.Range("E21").Value = '  -0.04%  '
$ws.Range("D22").Value = '2.103.23'
$ws.Range("E22").Value = '  -0.15%  '
$ws.Range("D23").Value = '''1.000'
$ws.Range("E23").Value = '  +0.02%  '
$ws.Range("D24").Value = '''7.403'
$ws.Range("E24").Value = '  -2.60%  '
$ws.Range("D25").Value = '''162.40'
$ws.Range("E25").Value = '  +0.70%  '
$ws.Range("D26").Value = '''8.963'
$ws.Range("E26").Value = '  -0.42%  '
$ws.Range("D27").Value = '''0.1447'
$ws.Range("E27").Value = '  -0.46%  '
$ws.Range("D28").Value = '''18.09'
$ws.Range("E28").Value = '  +0.12%  '
$ws.Range("D29").Value = '''1.957'
$ws.Range("E29").Value = '  -1.75%  '
$ws.Range("D30").Value = '''1.425'
$ws.Range("E30").Value = '  +2.05%  '
$ws.Range("D31").Value = '''1.485'
$ws.Range("E31").Value = '  -0.66%  '
$ws.Range("E32").Value = '  -3.65%  '
$ws.Range("E33").Value = '  +2.02%  '
$ws.Range("D34").Value = '''0.05215'
$ws.Range("E34").Value = '  +0.91%  '
$ws.Range("D35").Value = '''1.171'
$ws.Range("E35").Value = '  -1.30%  '
$ws.Range("D36").Value = '''0.7076'
$ws.Range("E36").Value = '  +0.69%  '
$ws.Range("D37").Value = '''0.9992'
$ws.Range("E37").Value = '  -4.18%  '
$ws.Range("D38").Value = '''2.672'
$ws.Range("E38").Value = '  +0.73%  '
$ws.Range("D39").Value = '''0.01849'
$ws.Range("E39").Value = '  -0.89%  '
$ws.Range("D40").Value = '''2.730'
$ws.Range("E40").Value = '  +1.97%  '
$ws.Range("D41").Value = '''0.9245'
$ws.Range("E41").Value = '  -1.57%  '
$ws.Range("D42").Value = '1.143.40'
$ws.Range("E42").Value = '  +6.15%  '
$ws.Range("D43").Value = '''0.4282'
$ws.Range("E43").Value = '  -0.09%  '
$ws.Range("D44").Value = '''5.871'
$ws.Range("E44").Value = '  -2.51%  '
$ws.Range("D45").Value = '''70.22'
$ws.Range("E45").Value = '  -0.09%  '
$ws.Range("D46").Value = '''0.9992'
$ws.Range("E46").Value = '  -0.02%  '
$ws.Range("D47").Value = '''102.79'
$ws.Range("E47").Value = '  +0.15%  '
$ws.Range("E48").Value = '  +1.95%  '
$ws.Range("D49").Value = '1.997.87'
$ws.Range("E49").Value = '  -0.22%  '
$ws.Range("D50").Value = '''9.206'
$ws.Range("E50").Value = '  +0.04%  '
$ws.Range("D51").Value = '''6.963'
$ws.Range("E51").Value = '  -1.20%  '
